# Insert a new weekly price-record row at row 52 (shifting existing rows
# 52..114 down to 53..115), then populate the new row with the latest
# observation for "Feria Lagunitas de Puerto Montt" / Ciboulette.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..114 down by one row.
$ws.Rows.Item(52).Insert()

# Fill in the newly inserted row 52 with the new weekly record.
$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 44467
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = 100112039
$ws.Cells.Item(52, 7).Value = "Ciboulette"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 240
$ws.Cells.Item(52, 11).Value = 3500
$ws.Cells.Item(52, 12).Value = 3500
$ws.Cells.Item(52, 13).Value = 3500
$ws.Cells.Item(52, 14).Value = "$/docena de atados"
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value = 1167
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"
